$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.033.38'
$ws.Range('E2').Value = '  -0.79%  '
$ws.Range('D3').Value = '2.583.13'
$ws.Range('E3').Value = '  -2.47%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '518.58'
$ws.Range('E5').Value = '  -0.63%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '142.67'
$ws.Range('E6').Value = '  -1.26%  '
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('E8').Value = '  -0.67%  '
$ws.Range('D9').Value = '2.599.38'
$ws.Range('E9').Value = '  -2.06%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.80'
$ws.Range('E10').Value = '  +1.42%  '
$ws.Range('E11').Value = '  -1.55%  '
$ws.Range('E12').Value = '  -3.47%  '
$ws.Range('E13').Value = '  -1.13%  '
$ws.Range('D14').Value = '3.039.49'
$ws.Range('E14').Value = '  -2.26%  '
$ws.Range('D15').Value = '58.021.60'
$ws.Range('E15').Value = '  -0.78%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '20.36'
$ws.Range('E16').Value = '  -2.40%  '
$ws.Range('E17').Value = '  -1.92%  '
$ws.Range('D18').Value = '2.571.01'
$ws.Range('E18').Value = '  -2.96%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '345.31'
$ws.Range('E19').Value = '  +2.01%  '
$ws.Range('E20').Value = '  -2.28%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '10.29'
$ws.Range('E21').Value = '  -1.76%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.38'
$ws.Range('E22').Value = '  +1.21%  '
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '66.06'
$ws.Range('E24').Value = '  +2.11%  '
$ws.Range('E25').Value = '  -1.23%  '
$ws.Range('E26').Value = '  -5.25%  '
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('D28').Value = '2.694.45'
$ws.Range('E28').Value = '  -2.47%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.00'
$ws.Range('E29').Value = '  -1.72%  '
$ws.Range('E30').Value = '  -6.81%  '
$ws.Range('E31').Value = '  -0.02%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.21'
$ws.Range('E32').Value = '  -6.33%  '
$ws.Range('E33').Value = '  -0.37%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '18.77'
$ws.Range('E34').Value = '  -0.50%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '149.57'
$ws.Range('E35').Value = '  -1.96%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.03'
$ws.Range('E36').Value = '  -2.63%  '
$ws.Range('E37').Value = '  -3.22%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.871'
$ws.Range('E38').Value = '  -4.04%  '
$ws.Range('B39').Value = 'Fetch.AI'
$ws.Range('C39').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.836'
$ws.Range('E39').Value = '  -2.68%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '35.95'
$ws.Range('E40').Value = '  -2.01%  '
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.54'
$ws.Range('E42').Value = '  -2.85%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.997'
$ws.Range('E43').Value = '  -0.29%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '273.74'
$ws.Range('E44').Value = '  +1.50%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '10.66'
$ws.Range('E45').Value = '  +0.23%  '
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.590'
$ws.Range('E46').Value = '  -3.00%  '
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0953'
$ws.Range('E47').Value = '  -1.65%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '18.83'
$ws.Range('E48').Value = '  -2.93%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0525'
$ws.Range('E49').Value = '  -2.30%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '4.66'
$ws.Range('E50').Value = '  -0.99%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '1.975.30'
$ws.Range('E51').Value = '  -3.35%  '
